# "Added totals column for cel/for/pot"
#
# Adds a Totals column (F) to Sheet1's three summary blocks (Energy
# Distribution, Type Distribution, Rarity Distribution), each new cell
# being =SUM(B<row>:E<row>) of the existing Core/Celerity/Fortitude/Potence
# columns. Also fixes the smart-quoted "4+" label in A9 to a plain "4+".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the curly-quoted label so it reads plainly.
$ws.Range("A9").Value = "4+"

# Pick up the existing column-B number formatting (style id 2 in the
# original file) for every new Totals cell by copying formats from the
# neighbouring Core column, then fill in the SUM formulas.

# Energy Distribution block (rows 4-10)
$ws.Range("B4:B10").Copy()
$ws.Range("F4:F10").PasteSpecial(-4122)

# Blank separator cell in row 12 (only formatted, no formula/value)
$ws.Range("B12").Copy()
$ws.Range("F12").PasteSpecial(-4122)

# Type Distribution block (rows 13-16)
$ws.Range("B13:B16").Copy()
$ws.Range("F13:F16").PasteSpecial(-4122)

# Rarity Distribution block (rows 19-23)
$ws.Range("B19:B23").Copy()
$ws.Range("F19:F23").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$totalRows = @(4,5,6,7,8,9,10,13,14,15,16,19,20,21,22,23)
foreach ($r in $totalRows) {
    $ws.Range("F$r").Formula = "=SUM(B" + $r + ":E" + $r + ")"
}

# Match the author's final selection / active cell.
$ws.Activate() | Out-Null
$ws.Range("F4").Select() | Out-Null
